# Update the "Marking" (B11) and "Total" (B12) scores, and the
# Correct/Total summary label in E12, on the "quiz" marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 40
$ws.Range("E12").Value = "40/140"
